# Spark_Part_List.xlsx update:
#  - Resistor R4 moved from the "100k" group to the "1k" group
#  - Document version date bumped
#  - Selection cursor moved to C29
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# 100k resistors: remove R4 from the reference designator list, qty 9 -> 8
$ws.Range("C25").Value = "R7, R10, R15, R16, R19, R38, R39, R40, R41, "
$ws.Range("A25").Value = 8

# Document version date
$ws.Range("C5").Value = "Document Version 01/06/2023"

# 1k resistors: add R4 to the reference designator list, qty 5 -> 6
$ws.Range("A22").Value = 6
$ws.Range("C22").Value = "R1, R4, R8, R24, R30, R37, "

# Update the active selection saved in the sheet view
$ws.Range("C29").Select()
